$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recommandations")

# Row 2: SUCRIVOIRE
$ws.Cells.Item(2,1).Value = "SUCRIVOIRE"
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = 5
$ws.Cells.Item(2,4).Value = 4920
$ws.Cells.Item(2,5).Value = 995
$ws.Cells.Item(2,6).Value = "🟡 Observer"
$ws.Cells.Item(2,7).Value = "➖ Neutre"

# Row 3: BRVM - SERVICES PUBLICS
$ws.Cells.Item(3,1).Value = "BRVM - SERVICES PUBLICS"
$ws.Cells.Item(3,2).Value = 0
$ws.Cells.Item(3,3).Value = 10
$ws.Cells.Item(3,4).Value = 4170.01
$ws.Cells.Item(3,5).Value = 112.02
$ws.Cells.Item(3,6).Value = "🟡 Observer"
$ws.Cells.Item(3,7).Value = "➖ Neutre"

# Row 4: SAFCA CI
$ws.Cells.Item(4,1).Value = "SAFCA CI"
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 5
$ws.Cells.Item(4,4).Value = 3470
$ws.Cells.Item(4,5).Value = 695
$ws.Cells.Item(4,6).Value = "🟡 Observer"
$ws.Cells.Item(4,7).Value = "➖ Neutre"

# Row 5: CFAO MOTORS CI
$ws.Cells.Item(5,1).Value = "CFAO MOTORS CI"
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 5
$ws.Cells.Item(5,4).Value = 3365
$ws.Cells.Item(5,5).Value = 675
$ws.Cells.Item(5,6).Value = "🟡 Observer"
$ws.Cells.Item(5,7).Value = "➖ Neutre"

# Row 6: BRVM - AUTRES SECTEURS
$ws.Cells.Item(6,1).Value = "BRVM - AUTRES SECTEURS"
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(6,3).Value = 5
$ws.Cells.Item(6,4).Value = 3297.09
$ws.Cells.Item(6,5).Value = 653.39
$ws.Cells.Item(6,6).Value = "🟡 Observer"
$ws.Cells.Item(6,7).Value = "➖ Neutre"

# Row 7: NEI-CEDA CI
$ws.Cells.Item(7,1).Value = "NEI-CEDA CI"
$ws.Cells.Item(7,2).Value = 0
$ws.Cells.Item(7,3).Value = 5
$ws.Cells.Item(7,4).Value = 2965
$ws.Cells.Item(7,5).Value = 595
$ws.Cells.Item(7,6).Value = "🟡 Observer"
$ws.Cells.Item(7,7).Value = "➖ Neutre"

# Row 8: UNIWAX CI
$ws.Cells.Item(8,1).Value = "UNIWAX CI"
$ws.Cells.Item(8,2).Value = 0
$ws.Cells.Item(8,3).Value = 5
$ws.Cells.Item(8,4).Value = 2910
$ws.Cells.Item(8,5).Value = 580
$ws.Cells.Item(8,6).Value = "🟡 Observer"
$ws.Cells.Item(8,7).Value = "➖ Neutre"

# Row 9: SETAO CI
$ws.Cells.Item(9,1).Value = "SETAO CI"
$ws.Cells.Item(9,2).Value = 0
$ws.Cells.Item(9,3).Value = 5
$ws.Cells.Item(9,4).Value = 2780
$ws.Cells.Item(9,5).Value = 555
$ws.Cells.Item(9,6).Value = "🟡 Observer"
$ws.Cells.Item(9,7).Value = "➖ Neutre"

# Row 10: AIR LIQUIDE CI
$ws.Cells.Item(10,1).Value = "AIR LIQUIDE CI"
$ws.Cells.Item(10,2).Value = 0
$ws.Cells.Item(10,3).Value = 5
$ws.Cells.Item(10,4).Value = 2650
$ws.Cells.Item(10,5).Value = 525
$ws.Cells.Item(10,6).Value = "🟡 Observer"
$ws.Cells.Item(10,7).Value = "➖ Neutre"

# Row 11: BRVM - DISTRIBUTION
$ws.Cells.Item(11,1).Value = "BRVM - DISTRIBUTION"
$ws.Cells.Item(11,2).Value = 0
$ws.Cells.Item(11,3).Value = 5
$ws.Cells.Item(11,4).Value = 1843.05
$ws.Cells.Item(11,5).Value = 367.04
$ws.Cells.Item(11,6).Value = "🟡 Observer"
$ws.Cells.Item(11,7).Value = "➖ Neutre"

# Row 12: BRVM - TRANSPORT
$ws.Cells.Item(12,1).Value = "BRVM - TRANSPORT"
$ws.Cells.Item(12,2).Value = 0
$ws.Cells.Item(12,3).Value = 5
$ws.Cells.Item(12,4).Value = 1741.56
$ws.Cells.Item(12,5).Value = 348.8
$ws.Cells.Item(12,6).Value = "🟡 Observer"
$ws.Cells.Item(12,7).Value = "➖ Neutre"

# Row 13: BRVM - AGRICULTURE
$ws.Cells.Item(13,1).Value = "BRVM - AGRICULTURE"
$ws.Cells.Item(13,2).Value = 0
$ws.Cells.Item(13,3).Value = 5
$ws.Cells.Item(13,4).Value = 1596.36
$ws.Cells.Item(13,5).Value = 308.95
$ws.Cells.Item(13,6).Value = "🟡 Observer"
$ws.Cells.Item(13,7).Value = "➖ Neutre"

# Row 14: BRVM - INDUSTRIE
$ws.Cells.Item(14,1).Value = "BRVM - INDUSTRIE"
$ws.Cells.Item(14,2).Value = 0
$ws.Cells.Item(14,3).Value = 5
$ws.Cells.Item(14,4).Value = 1003.91
$ws.Cells.Item(14,5).Value = 202.29
$ws.Cells.Item(14,6).Value = "🟡 Observer"
$ws.Cells.Item(14,7).Value = "➖ Neutre"

# Row 15: BRVM-PRINCIPAL
$ws.Cells.Item(15,1).Value = "BRVM-PRINCIPAL"
$ws.Cells.Item(15,2).Value = 0
$ws.Cells.Item(15,3).Value = 5
$ws.Cells.Item(15,4).Value = 878.35
$ws.Cells.Item(15,5).Value = 176.59
$ws.Cells.Item(15,6).Value = "🟡 Observer"
$ws.Cells.Item(15,7).Value = "➖ Neutre"

# Row 16: BRVM - CONSOMMATION DE BASE
$ws.Cells.Item(16,1).Value = "BRVM - CONSOMMATION DE BASE"
$ws.Cells.Item(16,2).Value = 0
$ws.Cells.Item(16,3).Value = 5
$ws.Cells.Item(16,4).Value = 872.15
$ws.Cells.Item(16,5).Value = 173.98
$ws.Cells.Item(16,6).Value = "🟡 Observer"
$ws.Cells.Item(16,7).Value = "➖ Neutre"

# Row 17: BRVM - INDUSTRIELS
$ws.Cells.Item(17,1).Value = "BRVM - INDUSTRIELS"
$ws.Cells.Item(17,2).Value = 0
$ws.Cells.Item(17,3).Value = 5
$ws.Cells.Item(17,4).Value = 663.35
$ws.Cells.Item(17,5).Value = 131.65
$ws.Cells.Item(17,6).Value = "🟡 Observer"
$ws.Cells.Item(17,7).Value = "➖ Neutre"

# Row 18: BRVM-PRESTIGE
$ws.Cells.Item(18,1).Value = "BRVM-PRESTIGE"
$ws.Cells.Item(18,2).Value = 0
$ws.Cells.Item(18,3).Value = 5
$ws.Cells.Item(18,4).Value = 650.62
$ws.Cells.Item(18,5).Value = 130.36
$ws.Cells.Item(18,6).Value = "🟡 Observer"
$ws.Cells.Item(18,7).Value = "➖ Neutre"

# Row 19: BRVM - FINANCES
$ws.Cells.Item(19,1).Value = "BRVM - FINANCES"
$ws.Cells.Item(19,2).Value = 0
$ws.Cells.Item(19,3).Value = 5
$ws.Cells.Item(19,4).Value = 607.78
$ws.Cells.Item(19,5).Value = 122.78
$ws.Cells.Item(19,6).Value = "🟡 Observer"
$ws.Cells.Item(19,7).Value = "➖ Neutre"

# Row 20: BRVM - SERVICES FINANCIERS
$ws.Cells.Item(20,1).Value = "BRVM - SERVICES FINANCIERS"
$ws.Cells.Item(20,2).Value = 0
$ws.Cells.Item(20,3).Value = 5
$ws.Cells.Item(20,4).Value = 597.33
$ws.Cells.Item(20,5).Value = 120.67
$ws.Cells.Item(20,6).Value = "🟡 Observer"
$ws.Cells.Item(20,7).Value = "➖ Neutre"

# Row 21: BRVM - ENERGIE
$ws.Cells.Item(21,1).Value = "BRVM - ENERGIE"
$ws.Cells.Item(21,2).Value = 0
$ws.Cells.Item(21,3).Value = 5
$ws.Cells.Item(21,4).Value = 546.48
$ws.Cells.Item(21,5).Value = 108.33
$ws.Cells.Item(21,6).Value = "🟡 Observer"
$ws.Cells.Item(21,7).Value = "➖ Neutre"

# Row 22: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws.Cells.Item(22,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws.Cells.Item(22,2).Value = 0
$ws.Cells.Item(22,3).Value = 5
$ws.Cells.Item(22,4).Value = 528.79
$ws.Cells.Item(22,5).Value = 105.5
$ws.Cells.Item(22,6).Value = "🟡 Observer"
$ws.Cells.Item(22,7).Value = "➖ Neutre"

# Row 23: BRVM - TELECOMMUNICATIONS
$ws.Cells.Item(23,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws.Cells.Item(23,2).Value = 0
$ws.Cells.Item(23,3).Value = 5
$ws.Cells.Item(23,4).Value = 469.87
$ws.Cells.Item(23,5).Value = 93.54
$ws.Cells.Item(23,6).Value = "🟡 Observer"
$ws.Cells.Item(23,7).Value = "➖ Neutre"

# Row 24: TRACTAFRIC MOTORS CI (PRSC)
$ws.Cells.Item(24,1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws.Cells.Item(24,2).Value = 3
$ws.Cells.Item(24,3).Value = 0
$ws.Cells.Item(24,4).Value = 17.5
$ws.Cells.Item(24,5).Value = 7.5
$ws.Cells.Item(24,6).Value = "🟢 Achat"
$ws.Cells.Item(24,7).Value = "✅ Renforcer"

# Row 25: BANK OF AFRICA ML (BOAM)
$ws.Cells.Item(25,1).Value = "BANK OF AFRICA ML (BOAM)"
$ws.Cells.Item(25,2).Value = 2
$ws.Cells.Item(25,3).Value = 0
$ws.Cells.Item(25,4).Value = 11.78
$ws.Cells.Item(25,5).Value = 4.99
$ws.Cells.Item(25,6).Value = "🟡 Observer"
$ws.Cells.Item(25,7).Value = "➖ Neutre"

# Row 26: BANK OF AFRICA NG (BOAN)
$ws.Cells.Item(26,1).Value = "BANK OF AFRICA NG (BOAN)"
$ws.Cells.Item(26,2).Value = 2
$ws.Cells.Item(26,3).Value = 0
$ws.Cells.Item(26,4).Value = 9.23
$ws.Cells.Item(26,5).Value = 3.4
$ws.Cells.Item(26,6).Value = "🟡 Observer"
$ws.Cells.Item(26,7).Value = "➖ Neutre"

# Row 27: SOCIETE IVOIRIENNE DE BANQUE  (SIBC)
$ws.Cells.Item(27,1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws.Cells.Item(27,2).Value = 1
$ws.Cells.Item(27,3).Value = 0
$ws.Cells.Item(27,4).Value = 5.59
$ws.Cells.Item(27,5).Value = 5.59
$ws.Cells.Item(27,6).Value = "🟡 Observer"
$ws.Cells.Item(27,7).Value = "➖ Neutre"

# Row 28: NSIA BANQUE COTE D'IVOIRE (NSBC)
$ws.Cells.Item(28,1).Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws.Cells.Item(28,2).Value = 1
$ws.Cells.Item(28,3).Value = 0
$ws.Cells.Item(28,4).Value = 3.3
$ws.Cells.Item(28,5).Value = 3.3
$ws.Cells.Item(28,6).Value = "🟡 Observer"
$ws.Cells.Item(28,7).Value = "➖ Neutre"

# Row 29: BANK OF AFRICA BN (BOAB)
$ws.Cells.Item(29,1).Value = "BANK OF AFRICA BN (BOAB)"
$ws.Cells.Item(29,2).Value = 1
$ws.Cells.Item(29,3).Value = 0
$ws.Cells.Item(29,4).Value = 3.09
$ws.Cells.Item(29,5).Value = 3.09
$ws.Cells.Item(29,6).Value = "🟡 Observer"
$ws.Cells.Item(29,7).Value = "➖ Neutre"

# Row 30: NEI-CEDA CI (NEIC)
$ws.Cells.Item(30,1).Value = "NEI-CEDA CI (NEIC)"
$ws.Cells.Item(30,2).Value = 1
$ws.Cells.Item(30,3).Value = 0
$ws.Cells.Item(30,4).Value = 2.59
$ws.Cells.Item(30,5).Value = 2.59
$ws.Cells.Item(30,6).Value = "🟡 Observer"
$ws.Cells.Item(30,7).Value = "➖ Neutre"

# Row 31: TOTALENERGIES MARKETING SN (TTLS)
$ws.Cells.Item(31,1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws.Cells.Item(31,2).Value = 2
$ws.Cells.Item(31,3).Value = 1
$ws.Cells.Item(31,4).Value = 1.41
$ws.Cells.Item(31,5).Value = -2.5
$ws.Cells.Item(31,6).Value = "🟡 Observer"
$ws.Cells.Item(31,7).Value = "👀 À surveiller"

# Row 32: BERNABE CI (BNBC)
$ws.Cells.Item(32,1).Value = "BERNABE CI (BNBC)"
$ws.Cells.Item(32,2).Value = 2
$ws.Cells.Item(32,3).Value = 2
$ws.Cells.Item(32,4).Value = 0.97
$ws.Cells.Item(32,5).Value = 4.74
$ws.Cells.Item(32,6).Value = "🟡 Observer"
$ws.Cells.Item(32,7).Value = "👀 À surveiller"

# Row 33: ECOBANK COTE D''IVOIRE (ECOC)
$ws.Cells.Item(33,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws.Cells.Item(33,2).Value = 1
$ws.Cells.Item(33,3).Value = 1
$ws.Cells.Item(33,4).Value = 0.28
$ws.Cells.Item(33,5).Value = 5.36
$ws.Cells.Item(33,6).Value = "🟡 Observer"
$ws.Cells.Item(33,7).Value = "👀 À surveiller"

# Row 34: SICABLE CI (CABC)
$ws.Cells.Item(34,1).Value = "SICABLE CI (CABC)"
$ws.Cells.Item(34,2).Value = 1
$ws.Cells.Item(34,3).Value = 1
$ws.Cells.Item(34,4).Value = 0.27
$ws.Cells.Item(34,5).Value = 5.31
$ws.Cells.Item(34,6).Value = "🟡 Observer"
$ws.Cells.Item(34,7).Value = "👀 À surveiller"

# Row 35: SONATEL SN (SNTS)
$ws.Cells.Item(35,1).Value = "SONATEL SN (SNTS)"
$ws.Cells.Item(35,2).Value = 1
$ws.Cells.Item(35,3).Value = 2
$ws.Cells.Item(35,4).Value = 0.25
$ws.Cells.Item(35,5).Value = -1.92
$ws.Cells.Item(35,6).Value = "🟡 Observer"
$ws.Cells.Item(35,7).Value = "👀 À surveiller"

# Row 36: TOTAL
$ws.Cells.Item(36,1).Value = "TOTAL"
$ws.Cells.Item(36,2).Value = 0
$ws.Cells.Item(36,3).Value = 5
$ws.Cells.Item(36,4).Value = 0
$ws.Cells.Item(36,5).Value = 0
$ws.Cells.Item(36,6).Value = "🟡 Observer"
$ws.Cells.Item(36,7).Value = "➖ Neutre"

# Row 37: SODE CI (SDCC)
$ws.Cells.Item(37,1).Value = "SODE CI (SDCC)"
$ws.Cells.Item(37,2).Value = 1
$ws.Cells.Item(37,3).Value = 1
$ws.Cells.Item(37,4).Value = -0.69
$ws.Cells.Item(37,5).Value = 3.36
$ws.Cells.Item(37,6).Value = "🟡 Observer"
$ws.Cells.Item(37,7).Value = "👀 À surveiller"

# Row 38: AFRICA GLOBAL LOGISTICS CI (SDSC)
$ws.Cells.Item(38,1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws.Cells.Item(38,2).Value = 1
$ws.Cells.Item(38,3).Value = 1
$ws.Cells.Item(38,4).Value = -0.71
$ws.Cells.Item(38,5).Value = -2.11
$ws.Cells.Item(38,6).Value = "🟡 Observer"
$ws.Cells.Item(38,7).Value = "👀 À surveiller"

# Row 39: VIVO ENERGY CI (SHEC)
$ws.Cells.Item(39,1).Value = "VIVO ENERGY CI (SHEC)"
$ws.Cells.Item(39,2).Value = 1
$ws.Cells.Item(39,3).Value = 2
$ws.Cells.Item(39,4).Value = -0.83
$ws.Cells.Item(39,5).Value = -2.35
$ws.Cells.Item(39,6).Value = "🟡 Observer"
$ws.Cells.Item(39,7).Value = "👀 À surveiller"

# Row 40: FILTISAC CI (FTSC)
$ws.Cells.Item(40,1).Value = "FILTISAC CI (FTSC)"
$ws.Cells.Item(40,2).Value = 1
$ws.Cells.Item(40,3).Value = 1
$ws.Cells.Item(40,4).Value = -0.91
$ws.Cells.Item(40,5).Value = -1.79
$ws.Cells.Item(40,6).Value = "🟡 Observer"
$ws.Cells.Item(40,7).Value = "👀 À surveiller"

# Row 41: BICI CI (BICC)
$ws.Cells.Item(41,1).Value = "BICI CI (BICC)"
$ws.Cells.Item(41,2).Value = 0
$ws.Cells.Item(41,3).Value = 1
$ws.Cells.Item(41,4).Value = -2.54
$ws.Cells.Item(41,5).Value = -2.54
$ws.Cells.Item(41,6).Value = "🟡 Observer"
$ws.Cells.Item(41,7).Value = "➖ Neutre"

# Row 42: SAFCA CI (SAFC)
$ws.Cells.Item(42,1).Value = "SAFCA CI (SAFC)"
$ws.Cells.Item(42,2).Value = 1
$ws.Cells.Item(42,3).Value = 1
$ws.Cells.Item(42,4).Value = -2.74
$ws.Cells.Item(42,5).Value = -6.47
$ws.Cells.Item(42,6).Value = "🟡 Observer"
$ws.Cells.Item(42,7).Value = "👀 À surveiller"

# Row 43: BANK OF AFRICA BF (BOABF)
$ws.Cells.Item(43,1).Value = "BANK OF AFRICA BF (BOABF)"
$ws.Cells.Item(43,2).Value = 0
$ws.Cells.Item(43,3).Value = 1
$ws.Cells.Item(43,4).Value = -2.86
$ws.Cells.Item(43,5).Value = -2.86
$ws.Cells.Item(43,6).Value = "🟡 Observer"
$ws.Cells.Item(43,7).Value = "➖ Neutre"

# Row 44: BANK OF AFRICA SENEGAL (BOAS)
$ws.Cells.Item(44,1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws.Cells.Item(44,2).Value = 0
$ws.Cells.Item(44,3).Value = 1
$ws.Cells.Item(44,4).Value = -3.37
$ws.Cells.Item(44,5).Value = -3.37
$ws.Cells.Item(44,6).Value = "🟡 Observer"
$ws.Cells.Item(44,7).Value = "➖ Neutre"

# Row 45: ORANGE COTE D'IVOIRE (ORAC)
$ws.Cells.Item(45,1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws.Cells.Item(45,2).Value = 0
$ws.Cells.Item(45,3).Value = 1
$ws.Cells.Item(45,4).Value = -3.38
$ws.Cells.Item(45,5).Value = -3.38
$ws.Cells.Item(45,6).Value = "🟡 Observer"
$ws.Cells.Item(45,7).Value = "➖ Neutre"

# Row 46: UNIWAX CI (UNXC)
$ws.Cells.Item(46,1).Value = "UNIWAX CI (UNXC)"
$ws.Cells.Item(46,2).Value = 1
$ws.Cells.Item(46,3).Value = 2
$ws.Cells.Item(46,4).Value = -3.84
$ws.Cells.Item(46,5).Value = -4.35
$ws.Cells.Item(46,6).Value = "🟡 Observer"
$ws.Cells.Item(46,7).Value = "👀 À surveiller"

# Row 47: SOLIBRA CI (SLBC)
$ws.Cells.Item(47,1).Value = "SOLIBRA CI (SLBC)"
$ws.Cells.Item(47,2).Value = 1
$ws.Cells.Item(47,3).Value = 2
$ws.Cells.Item(47,4).Value = -4.85
$ws.Cells.Item(47,5).Value = 7.5
$ws.Cells.Item(47,6).Value = "🟡 Observer"
$ws.Cells.Item(47,7).Value = "👀 À surveiller"

# Row 48: CORIS BANK INTERNATIONAL (CBIBF)
$ws.Cells.Item(48,1).Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws.Cells.Item(48,2).Value = 0
$ws.Cells.Item(48,3).Value = 1
$ws.Cells.Item(48,4).Value = -5.69
$ws.Cells.Item(48,5).Value = -5.69
$ws.Cells.Item(48,6).Value = "🟡 Observer"
$ws.Cells.Item(48,7).Value = "➖ Neutre"

# Row 49: ECOBANK TRANS. INCORP. TG (ETIT)
$ws.Cells.Item(49,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws.Cells.Item(49,2).Value = 0
$ws.Cells.Item(49,3).Value = 1
$ws.Cells.Item(49,4).Value = -5.88
$ws.Cells.Item(49,5).Value = -5.88
$ws.Cells.Item(49,6).Value = "🟡 Observer"
$ws.Cells.Item(49,7).Value = "➖ Neutre"

# Row 50: SETAO CI (STAC)
$ws.Cells.Item(50,1).Value = "SETAO CI (STAC)"
$ws.Cells.Item(50,2).Value = 0
$ws.Cells.Item(50,3).Value = 1
$ws.Cells.Item(50,4).Value = -6.09
$ws.Cells.Item(50,5).Value = -6.09
$ws.Cells.Item(50,6).Value = "🟡 Observer"
$ws.Cells.Item(50,7).Value = "➖ Neutre"

# Row 51: SAPH CI (SPHC)
$ws.Cells.Item(51,1).Value = "SAPH CI (SPHC)"
$ws.Cells.Item(51,2).Value = 0
$ws.Cells.Item(51,3).Value = 1
$ws.Cells.Item(51,4).Value = -7.46
$ws.Cells.Item(51,5).Value = -7.46
$ws.Cells.Item(51,6).Value = "🟡 Observer"
$ws.Cells.Item(51,7).Value = "➖ Neutre"

$ws2 = $wb.Worksheets.Item("Top_YTD")
$ws2.Cells.Item(2,2).Value = 159772762.21
$ws2.Cells.Item(3,2).Value = 14965258.03
$ws2.Cells.Item(4,2).Value = 3154460
$ws2.Cells.Item(5,2).Value = 2759640.49
$ws2.Cells.Item(6,2).Value = 2525414.13
$ws2.Cells.Item(7,2).Value = 1597843.31
$ws2.Cells.Item(8,2).Value = 1474979.12
$ws2.Cells.Item(9,2).Value = 1213232.48
$ws2.Cells.Item(10,2).Value = 992274.03
$ws2.Cells.Item(11,2).Value = 225841.42
